# Final - 5th May 2025
# Adds a new "HoverIcon" worksheet (with the hover-icon tooltip copy for the
# "round trip" flag) positioned between "PicklistValues" and "Contact", and
# makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Insert the new worksheet immediately before "Contact" -> lands right after
# "PicklistValues", matching the target sheet order:
#   Users, PicklistValues, HoverIcon, Contact, AddOpportunity, AddContact
$ws = $wb.Worksheets.Add($wb.Worksheets.Item("Contact"))
$ws.Name = "HoverIcon"

# Column A is wide enough to hold the wrapped explanatory paragraph.
$ws.Columns.Item(1).ColumnWidth = 65.1

# Row 1: bold header/title cell.
$ws.Range("A1").Value = "Hover Icon Text"
$ws.Range("A1").Font.Bold = $true

# Row 2: the wrapped tooltip description text.
$ws.Range("A2").Value = "An engagement is typically considered a potential round trip if it is acquired by a sponsor (subject is a potential round trip) or by a sponsor-backed operating company (buyer is a potential round trip). Note ""sponsor"" includes firms tagged as Private Equity Group, Hedge Fund, or Family Office."
$ws.Range("A2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 57.6

# Match the saved selection/active state on the new sheet; since this sheet
# was just added it also becomes the active (selected) tab of the workbook.
$ws.Range("C6").Select() | Out-Null
